# Auto-generated edit script applying the Coeurl_Profits.xlsx diff
# Updates static numeric values across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1388.4546
$ws.Range("I18").Value = 1388.4546
$ws.Range("K18").Value = 1388.4546
$ws.Range("M18").Value = -1104.4546
$ws.Range("H28").Value = 2105.0557
$ws.Range("I28").Value = 2084.0833
$ws.Range("K28").Value = 2084.0833
$ws.Range("M28").Value = -1599.0833
$ws.Range("H43").Value = 980.5714
$ws.Range("J43").Value = 993.3333
$ws.Range("L43").Value = 993.3333
$ws.Range("N43").Value = -1131.3333
$ws.Range("H103").Value = 697
$ws.Range("J103").Value = 697
$ws.Range("L103").Value = 2091
$ws.Range("N103").Value = -3263
$ws.Range("H112").Value = 73225.92999999999
$ws.Range("I112").Value = 1419.6666
$ws.Range("K112").Value = 4258.9998
$ws.Range("M112").Value = -3150.9998
$ws.Range("H132").Value = 1836.5
$ws.Range("I132").Value = 1537.6154
$ws.Range("J132").Value = 2896.182
$ws.Range("K132").Value = 4612.8462
$ws.Range("L132").Value = 8688.545999999998
$ws.Range("M132").Value = -2082.8462
$ws.Range("N132").Value = -13748.546
$ws.Range("H134").Value = 310000
$ws.Range("J134").Value = 310000
$ws.Range("L134").Value = 310000
$ws.Range("N134").Value = -320140

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2795.1792
$ws.Range("I61").Value = 2538.2886
$ws.Range("K61").Value = 2538.2886
$ws.Range("M61").Value = -2326.2886
$ws.Range("H74").Value = 5191.7827
$ws.Range("I74").Value = 1190.2106
$ws.Range("K74").Value = 1190.2106
$ws.Range("M74").Value = -316.2106000000001
$ws.Range("H77").Value = 5191.7827
$ws.Range("I77").Value = 1190.2106
$ws.Range("K77").Value = 5951.053000000001
$ws.Range("M77").Value = -1583.053000000001
$ws.Range("H102").Value = 3484.1
$ws.Range("I102").Value = 3010.7646
$ws.Range("K102").Value = 3010.7646
$ws.Range("M102").Value = -1388.7646
$ws.Range("H122").Value = 3058.2
$ws.Range("I122").Value = 2998.3044
$ws.Range("J122").Value = 3747
$ws.Range("K122").Value = 8994.913199999999
$ws.Range("L122").Value = 11241
$ws.Range("M122").Value = -6544.913199999999
$ws.Range("N122").Value = -16141
$ws.Range("H132").Value = 2681.0754
$ws.Range("I132").Value = 2572.525
$ws.Range("K132").Value = 7717.575000000001
$ws.Range("M132").Value = -5187.575000000001
$ws.Range("H136").Value = 2795.1792
$ws.Range("I136").Value = 2538.2886
$ws.Range("K136").Value = 7614.8658
$ws.Range("M136").Value = -5064.8658

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H134").Value = 1696.2456
$ws.Range("I134").Value = 1693.36
$ws.Range("J134").Value = 1716.8572
$ws.Range("K134").Value = 5080.08
$ws.Range("L134").Value = 5150.571599999999
$ws.Range("M134").Value = -2545.08
$ws.Range("N134").Value = -10220.5716

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 731705.7
$ws.Range("I3").Value = 1250585.2
$ws.Range("J3").Value = 39866.332
$ws.Range("K3").Value = 1250585.2
$ws.Range("L3").Value = 39866.332
$ws.Range("M3").Value = -1250472.2
$ws.Range("N3").Value = -40092.332
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -21250
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H60").Value = 11199.4
$ws.Range("I60").Value = 5333
$ws.Range("K60").Value = 5333
$ws.Range("M60").Value = -4822
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H132").Value = 3087.3967
$ws.Range("I132").Value = 2894.1091
$ws.Range("J132").Value = 4416.25
$ws.Range("K132").Value = 8682.327300000001
$ws.Range("L132").Value = 13248.75
$ws.Range("M132").Value = -6152.327300000001
$ws.Range("N132").Value = -18308.75
$ws.Range("H137").Value = 95000
$ws.Range("J137").Value = 95000
$ws.Range("L137").Value = 95000
$ws.Range("N137").Value = -105200

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1535.75
$ws.Range("J60").Value = 188.33333
$ws.Range("L60").Value = 564.99999
$ws.Range("N60").Value = -1066.99999
$ws.Range("H137").Value = 3215.389
$ws.Range("J137").Value = 4747.6665
$ws.Range("L137").Value = 14242.9995
$ws.Range("N137").Value = -24442.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H113").Value = 2104.3684
$ws.Range("I113").Value = 1918
$ws.Range("J113").Value = 2803.25
$ws.Range("K113").Value = 1918
$ws.Range("L113").Value = 2803.25
$ws.Range("M113").Value = 252
$ws.Range("N113").Value = -7143.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3111.6667
$ws.Range("I68").Value = 2488.889
$ws.Range("K68").Value = 2488.889
$ws.Range("M68").Value = -1739.889
$ws.Range("H71").Value = 3111.6667
$ws.Range("I71").Value = 2488.889
$ws.Range("K71").Value = 12444.445
$ws.Range("M71").Value = -8700.445
$ws.Range("H132").Value = 5457.136
$ws.Range("I132").Value = 4863
$ws.Range("J132").Value = 6170.1
$ws.Range("K132").Value = 14589
$ws.Range("L132").Value = 18510.3
$ws.Range("M132").Value = -12059
$ws.Range("N132").Value = -23570.3

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2097.925
$ws.Range("I122").Value = 2051.5757
$ws.Range("K122").Value = 6154.7271
$ws.Range("M122").Value = -3704.7271
$ws.Range("H132").Value = 1486.5358
$ws.Range("I132").Value = 1580.1666
$ws.Range("J132").Value = 924.75
$ws.Range("K132").Value = 4740.4998
$ws.Range("L132").Value = 2774.25
$ws.Range("M132").Value = -2210.4998
$ws.Range("N132").Value = -7834.25
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
$ws.Range("H136").Value = 2051.093
$ws.Range("I136").Value = 1550.2333
$ws.Range("J136").Value = 3206.923
$ws.Range("K136").Value = 4650.699900000001
$ws.Range("L136").Value = 9620.769
$ws.Range("M136").Value = -2100.699900000001
$ws.Range("N136").Value = -14720.769
$ws.Range("H137").Value = 100664.336
$ws.Range("J137").Value = 100664.336
$ws.Range("L137").Value = 100664.336
$ws.Range("N137").Value = -110864.336
$ws.Range("H139").Value = 81388
$ws.Range("I139").Value = 62777
$ws.Range("J139").Value = 118610
$ws.Range("K139").Value = 62777
$ws.Range("L139").Value = 118610
$ws.Range("M139").Value = -57637
$ws.Range("N139").Value = -128890
$ws.Range("H140").Value = 69755.8
$ws.Range("J140").Value = 69755.8
$ws.Range("L140").Value = 69755.8
$ws.Range("N140").Value = -80115.8
$ws.Range("H141").Value = 250000
$ws.Range("J141").Value = 250000
$ws.Range("L141").Value = 250000
$ws.Range("N141").Value = -260360
